# Updates the cryptocurrency price/volume table on the active worksheet
# to the latest scraped values (GitHub Actions scheduled refresh).
#
# Price values in column D are written with a leading apostrophe where the
# new text would otherwise be auto-recognized by Excel as a number (e.g.
# "5.28" or "0.0000237"), so that the cell keeps its original plain-text
# representation (matching values that use "." as both a thousands and a
# decimal separator, like "61.827.97", remain text on their own).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.827.97"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "3.054.26"
$ws.Range("E3").Value = "  -3.84%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'587.19"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").Value = "'130.46"
$ws.Range("E6").Value = "  -3.42%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.055.62"
$ws.Range("E8").Value = "  -3.76%  "
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("D11").Value = "'5.28"
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("D12").Value = "'0.442"
$ws.Range("E12").Value = "  -2.56%  "
$ws.Range("D13").Value = "'0.0000237"
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").Value = "'33.73"
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("D15").Value = "'0.120"
$ws.Range("E15").Value = "  +0.66%  "
$ws.Range("D16").Value = "3.557.54"
$ws.Range("E16").Value = "  -3.79%  "
$ws.Range("D17").Value = "61.870.48"
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("D18").Value = "3.049.83"
$ws.Range("E18").Value = "  -3.92%  "
$ws.Range("D19").Value = "'6.40"
$ws.Range("E19").Value = "  -2.11%  "
$ws.Range("D20").Value = "'449.83"
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("D21").Value = "'13.58"
$ws.Range("E21").Value = "  -2.36%  "
$ws.Range("D22").Value = "'0.675"
$ws.Range("E22").Value = "  -4.09%  "
$ws.Range("E23").Value = "  -3.25%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'81.28"
$ws.Range("E24").Value = "  -2.62%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "'12.90"
$ws.Range("E25").Value = "  -3.48%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'0.998"
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("E28").Value = "  -3.88%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").Value = "'7.44"
$ws.Range("E30").Value = "  -3.69%  "
$ws.Range("D31").Value = "'6.45"
$ws.Range("E31").Value = "  -6.15%  "
$ws.Range("D32").Value = "'25.97"
$ws.Range("E32").Value = "  -5.04%  "
$ws.Range("D33").Value = "'0.0978"
$ws.Range("E33").Value = "  -6.10%  "
$ws.Range("D34").Value = "'2.33"
$ws.Range("E34").Value = "  -2.21%  "
$ws.Range("D35").Value = "'0.978"
$ws.Range("E35").Value = "  -5.18%  "
$ws.Range("D36").Value = "'5.74"
$ws.Range("E36").Value = "  -2.53%  "
$ws.Range("D37").Value = "'50.45"
$ws.Range("E37").Value = "  -1.04%  "
$ws.Range("D38").Value = "0.0₃0697"
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("D39").Value = "'0.0377"
$ws.Range("E39").Value = "  -1.91%  "
$ws.Range("E40").Value = "  -0.51%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.109"
$ws.Range("E41").Value = "  -2.16%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "'380.72"
$ws.Range("E42").Value = "  -6.35%  "
$ws.Range("E43").Value = "  -5.78%  "
$ws.Range("D44").Value = "2.694.34"
$ws.Range("E44").Value = "  -6.16%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "'124.23"
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("D47").Value = "'0.241"
$ws.Range("E47").Value = "  -3.16%  "
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").Value = "'2.03"
$ws.Range("E48").Value = "  -4.73%  "
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").Value = "'34.11"
$ws.Range("E49").Value = "  -6.09%  "
$ws.Range("E50").Value = "  -1.71%  "
$ws.Range("D51").Value = "'24.09"
$ws.Range("E51").Value = "  -4.96%  "
